# Fix bug compare two objects
# - Replace the "Sheet1" tab content with a new 4-column comparison table
# - Move "Sheet1" to be the first tab (in front of "Trang tính1")

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Sheet1")
$otherSheet = $wb.Worksheets.Item("Trang tính1")

# --- Grab donor-formatted cells from $otherSheet before we touch anything,
# so we can stamp the exact same cell styles (s=1 / s=2 / s=3) onto the
# newly-created cells without creating brand-new style table entries.
$otherSheet.Range("A3").Copy()
$dataSheet.Range("C1").PasteSpecial(-4122)   # xlPasteFormats -> reuse style s=1 (same as A1/B1/A2)
$otherSheet.Range("B3").Copy()
$dataSheet.Range("B3").PasteSpecial(-4122)   # xlPasteFormats -> reuse style s=3
$otherSheet.Range("A3").Copy()
$dataSheet.Range("A3").PasteSpecial(-4122)   # xlPasteFormats -> reuse style s=1
$otherSheet.Range("A3").Copy()
$dataSheet.Range("C2").PasteSpecial(-4122)   # xlPasteFormats -> reuse style s=1
$otherSheet.Range("A3").Copy()
$dataSheet.Range("D1").PasteSpecial(-4122)   # xlPasteFormats -> reuse style s=1
$otherSheet.Range("A3").Copy()
$dataSheet.Range("C3").PasteSpecial(-4122)   # xlPasteFormats -> reuse style s=1

# --- Header row
$dataSheet.Range("A1").Value = "Facebook"
$dataSheet.Range("B1").Value = "user Name"
$dataSheet.Range("C1").Value = "password"
$dataSheet.Range("D1").Value = "PlaceHolder"

# --- Row 2
$dataSheet.Range("A2").Value = 1
$dataSheet.Range("B2").Value = "j"
$dataSheet.Range("C2").Value = "k"
$dataSheet.Range("D2").Value = 0

# --- Row 3
$dataSheet.Range("A3").Value = 3
$dataSheet.Range("B3").Value = "znpmdrz_goldmanberg_1473307997@tfbnw.net"
$dataSheet.Range("C3").Value = "ddd"
$dataSheet.Range("D3").Value = 0

# --- Column widths (best-fit-like): A, B, C get explicit widths; D stays default
$dataSheet.Columns.Item(1).ColumnWidth = 7.833333333333334
$dataSheet.Columns.Item(2).ColumnWidth = 39.66666666666667
$dataSheet.Columns.Item(3).ColumnWidth = 7.5

# --- Selection: whole column D selected, active cell D1
$dataSheet.Range("D1:D1048576").Select()
$dataSheet.Range("D1").Activate()

# --- Move "Sheet1" to the front (position 1) and keep it the active tab
$dataSheet.Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item(1).Activate()
